# [Update] Implementing dropdown feature
#
# Adds a new "Dropdown" column (E) to the Response/API mapping sheet.
# Column E holds the name of the dropdown-driven parameter associated
# with a handful of existing rows (end_use / equipment), and the
# previously-too-wide A/B columns are narrowed to make room for it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Dropdown" column header + values ------------------------------
$ws.Range("E1").Value = "Dropdown"
$ws.Range("E7").Value = "end_use"
$ws.Range("E8").Value = "equipment"
$ws.Range("E14").Value = "equipment"

# --- Column widths: narrow A/B, size the new E column --------------------
$ws.Columns.Item(1).ColumnWidth = 93.66666666666667   # A: ~94.43 chars
$ws.Columns.Item(2).ColumnWidth = 45.5                 # B: ~46.29 chars
$ws.Columns.Item(5).ColumnWidth = 19.833333333333332   # E: ~20.71 chars

# --- View state: scroll so column B is the leftmost visible column, and
#     leave the selection on the newly-edited E14 cell ------------------
$ws.Activate()
$ws.Range("E14").Select()
$excel.ActiveWindow.ScrollColumn = 2
